# Insert a new register row ("brd_id") above the existing "cr_dest_port"
# row (current row 26), shifting cr_dest_port and everything below it down
# by one row. Then update cr_dest_port's offset_from_msb (column D) to
# account for the 4 bits now used by brd_id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 26 (existing rows 26+ shift down to 27+).
$ws.Rows("26:26").Insert()

# Populate the new row 26 with the brd_id register definition.
$ws.Range("A26").Value = "brd_id"
$ws.Range("B26").Value = 4
$ws.Range("C26").Value = "cosmic_ray_eth_control"
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 32
$ws.Range("F26").Value = "readwrite"
$ws.Range("G26").Value = "Number to identify each SNAP2 board"

# cr_dest_port (now on row 27) now starts 4 bits later within the same
# mainregister, since brd_id occupies bits 3-6.
$ws.Range("D27").Value = 7

# Match the author's final on-screen selection.
$ws.Range("D27").Select() | Out-Null
